$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "60.472.91"
Set-TextValue "E2" "  +1.29%  "
Set-TextValue "D3" "2.584.89"
Set-TextValue "E3" "  +0.12%  "
Set-TextValue "E4" "  +0.23%  "
Set-TextValue "D5" "506.71"
Set-TextValue "E5" "  +0.64%  "
Set-TextValue "D6" "153.71"
Set-TextValue "E6" "  -1.53%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  +0.11%  "
Set-TextValue "D8" "0.583"
Set-TextValue "E8" "  -6.86%  "
Set-TextValue "D9" "2.591.83"
Set-TextValue "E9" "  +0.43%  "
Set-TextValue "D10" "6.57"
Set-TextValue "E10" "  +7.05%  "
Set-TextValue "E11" "  +0.97%  "
Set-TextValue "D12" "0.347"
Set-TextValue "E12" "  +2.28%  "
Set-TextValue "E13" "  +0.81%  "
Set-TextValue "D14" "3.041.97"
Set-TextValue "E14" "  +1.17%  "
Set-TextValue "D15" "60.499.68"
Set-TextValue "E15" "  +1.65%  "
Set-TextValue "D16" "21.49"
Set-TextValue "E16" "  -0.89%  "
Set-TextValue "E17" "  +2.06%  "
Set-TextValue "D18" "2.593.88"
Set-TextValue "E18" "  +0.77%  "
Set-TextValue "E20" "  +4.21%  "
Set-TextValue "D21" "10.43"
Set-TextValue "E21" "  +1.02%  "
Set-TextValue "E22" "  +1.58%  "
Set-TextValue "D23" "0.998"
Set-TextValue "E23" "  -1.13%  "
Set-TextValue "D24" "60.00"
Set-TextValue "E24" "  +0.32%  "
Set-TextValue "D25" "0.419"
Set-TextValue "E25" "  +1.10%  "
Set-TextValue "E26" "  +0.64%  "
Set-TextValue "D27" "0.996"
Set-TextValue "E27" "  +0.25%  "
Set-TextValue "D28" "0.0₃0845"
Set-TextValue "E28" "  +2.79%  "
Set-TextValue "D29" "7.35"
Set-TextValue "E29" "  -0.49%  "
Set-TextValue "E30" "  -0.03%  "
Set-TextValue "D31" "19.34"
Set-TextValue "E31" "  +0.26%  "
Set-TextValue "D32" "154.04"
Set-TextValue "E32" "  -2.12%  "
Set-TextValue "E33" "  -0.43%  "
Set-TextValue "D34" "5.69"
Set-TextValue "E34" "  +3.44%  "
Set-TextValue "E36" "  +0.18%  "
Set-TextValue "D37" "0.860"
Set-TextValue "E37" "  +10.96%  "
Set-TextValue "E38" "  +1.17%  "
Set-TextValue "E39" "  +0.59%  "
Set-TextValue "E40" "  +1.37%  "
Set-TextValue "E41" "  +2.24%  "
Set-TextValue "D42" "296.77"
Set-TextValue "E42" "  +1.91%  "
Set-TextValue "E43" "  -1.63%  "
Set-TextValue "D44" "0.613"
Set-TextValue "E44" "  -1.84%  "
Set-TextValue "E45" "  -0.23%  "
Set-TextValue "D46" "0.0556"
Set-TextValue "E46" "  -1.04%  "
Set-TextValue "D47" "19.77"
Set-TextValue "E47" "  +3.28%  "
Set-TextValue "E48" "  +1.14%  "
Set-TextValue "E49" "  -1.01%  "
Set-TextValue "D51" "2.002.10"
Set-TextValue "E51" "  +0.11%  "
